$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: status ("In Translation" -> "Ready for handoff") for both locale
# columns, plus the refreshed "Latest HO Xliff Generate Date" timestamp.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-12 15:11:40"

# zh-cn detail sheet: status + its own handoff datetime.
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-12 15:11:32"

# de-de detail sheet: status (its handoff datetime string matches the Overview
# sheet's refreshed generate-date value).
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-12 15:11:40"

# The new status text ("Ready for handoff") is wider than the old one, so the
# status columns grew to fit it.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
